# Libro1.xlsx edit: update a few cell values on "Hoja1" and switch the
# active sheet / selection from Hoja2!D8 to Hoja1!F8 (Hoja1 becomes the
# selected tab instead of Hoja2).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Hoja1")

# Update cell contents (order matters for shared-string table layout:
# F7 first, then E7, then D7, so the three new unique strings land in the
# same order as the target workbook).
$ws1.Range("F7").Value = "c/\"
$ws1.Range("E7").Value = 'a"\'
$ws1.Range("D7").Value = "s/"
$ws1.Range("F8").Value = "nn"

# Make Hoja1 the active sheet/tab and move its selection to F8.
$ws1.Activate()
$ws1.Range("F8").Select()
